$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): U1 / V1 get new labels, W1 is a new column header
$ws.Range("U1").Value = "param_E_pv3_solar"
$ws.Range("V1").Value = "param_P_to_charging_station1"
$ws.Range("W1").Value = "param_P_to_charging_station2"

# W1 is a brand new header cell - give it the same formatting (bold, border,
# centered) as the rest of the header row by copying the format from V1
$ws.Range("V1").Copy()
$ws.Range("W1").PasteSpecial(-4122)

# Data for columns U, V, W across rows 2-17
$data = @{
    2  = @{ U = 0.12; V = 0;                 W = 0 }
    3  = @{ U = 0.12; V = 48.4227638888889;  W = 0 }
    4  = @{ U = 0.12; V = 107.8514861111111; W = 16.21400000000001 }
    5  = @{ U = 0.12; V = 104.8002083333334; W = 139.1058472222223 }
    6  = @{ U = 0.12; V = 0;                 W = 0 }
    7  = @{ U = 0.12; V = 0;                 W = 0 }
    8  = @{ U = 0.12; V = 87.94500000000002; W = 44.99000000000002 }
    9  = @{ U = 0.12; V = 161.8688194444445; W = 113.5164861111111 }
    10 = @{ U = 0.12; V = 0;                 W = 0 }
    11 = @{ U = 0.12; V = 0;                 W = 47.85000000000002 }
    12 = @{ U = 0.12; V = 0;                 W = 109.395 }
    13 = @{ U = 0.12; V = 59.89148611111113; W = 49.33500000000002 }
    14 = @{ U = 0.12; V = 54.11648611111113; W = 0 }
    15 = @{ U = 0.12; V = 0;                 W = 0 }
    16 = @{ U = 0.12; V = 71.51084722222224; W = 0 }
    17 = @{ U = 0.12; V = 30.25000000000002; W = 0 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("U$row").Value = $vals.U
    $ws.Range("V$row").Value = $vals.V
    $ws.Range("W$row").Value = $vals.W
}
